# uml/cd/associationAsAttributes: correct association label direction
#
# Four small "multiplicity" triangles (square-headed association-label
# markers) on the single slide get their rotation/flip/position fixed so
# that the label direction is consistent. The triangle that used to sit
# in the 2nd spot keeps its original geometry "slot" being vacated while
# the rest shift down one, and a fresh placement is introduced for what
# used to be the last triangle in the chain.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.Id -eq $id) { return $cand }
    }
    return $null
}

# The four "Isosceles Triangle" association-multiplicity markers affected
# by this fix (identified by their current/old shape Id).
$triA = Get-ShapeById $s 51   # "Isosceles Triangle 50"
$triB = Get-ShapeById $s 52   # "Isosceles Triangle 51" -> becomes "...52"
$triC = Get-ShapeById $s 53   # "Isosceles Triangle 52" -> becomes "...53"
$triD = Get-ShapeById $s 54   # "Isosceles Triangle 53" -> becomes "...54"

# --- triA: rotate from 90 deg to 270 deg, flip horizontally, reposition ---
$triA.Rotation = 270
$triA.Flip(0)   # msoFlipHorizontal
$triA.Left = 164.07768357164278
$triA.Top = 118.75515530353456

# --- triB: same rotation/flip fix, takes over triC's old position, renamed ---
$triB.Rotation = 270
$triB.Flip(0)   # msoFlipHorizontal
$triB.Left = 201.43853826928327
$triB.Top = 166.55177424513448
$triB.Name = "Isosceles Triangle 52"

# --- triC: already rotated/flipped correctly, takes over triD's old position, renamed ---
$triC.Left = 621.422974420983
$triC.Top = 166.63838269631694
$triC.Name = "Isosceles Triangle 53"

# --- triD: already rotated/flipped correctly, moves to a brand-new position, renamed ---
$triD.Left = 581.9769901690145
$triD.Top = 120.08169359762837
$triD.Name = "Isosceles Triangle 54"
